# Modified templates to display posts from database
#
# Appends a block of new paragraphs (cryptoparty.pl / hackerspace.pl /
# girls code fun / py ladies / py ladies warsaw, with blank separator
# paragraphs) after the existing content, right before the section break.

$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Add-RawParagraph {
    # Appends a brand-new paragraph at the end of the document whose
    # contents are exactly $innerXml (the w:p children). Inserting a
    # plain paragraph mark first gives InsertXML a fresh, empty landing
    # paragraph to replace -- calling InsertXML directly on the existing
    # last paragraph would overwrite *its* content instead of adding a
    # new one.
    param([string]$innerXml = "")

    $d.Paragraphs.Last.Range.InsertParagraphAfter()
    $target = $d.Paragraphs.Last
    $target.Range.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

# blank separator paragraph
Add-RawParagraph

# cryptoparty.pl
Add-RawParagraph '<w:r><w:t>cryptoparty.pl</w:t></w:r>'

# hackerspace.pl
Add-RawParagraph '<w:r><w:t>hackerspace.pl</w:t></w:r>'

# blank separator paragraph
Add-RawParagraph

# girls code  fun
Add-RawParagraph ('<w:r><w:t xml:space="preserve">girls </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>code</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve">  </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>fun</w:t></w:r><w:proofErr w:type="spellEnd"/>')

# blank separator paragraph
Add-RawParagraph

# py ladies
Add-RawParagraph ('<w:proofErr w:type="spellStart"/><w:r><w:t>py</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>ladies</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>')

# py ladies warsaw
Add-RawParagraph ('<w:proofErr w:type="spellStart"/><w:r><w:t>py</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>ladies</w:t></w:r><w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/><w:r><w:t>warsaw</w:t></w:r><w:proofErr w:type="spellEnd"/>')

# trailing blank separator paragraph
Add-RawParagraph

Write-Host "Final paragraph count:" $d.Paragraphs.Count
